$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 31250536
$ws.Range("I19").Value = 233.57143
$ws.Range("J19").Value = 55556330
$ws.Range("K19").Value = 233.57143
$ws.Range("L19").Value = 55556330
$ws.Range("M19").Value = -58.57142999999999
$ws.Range("N19").Value = -55556680
$ws.Range("H70").Value = 1172.8611
$ws.Range("J70").Value = 833
$ws.Range("L70").Value = 2499
$ws.Range("N70").Value = -3039
$ws.Range("H73").Value = 1172.8611
$ws.Range("J73").Value = 833
$ws.Range("L73").Value = 2499
$ws.Range("N73").Value = -4371
$ws.Range("H86").Value = 71851.82000000001
$ws.Range("I86").Value = 93321.62
$ws.Range("J86").Value = 2075
$ws.Range("K86").Value = 93321.62
$ws.Range("L86").Value = 2075
$ws.Range("M86").Value = -92198.62
$ws.Range("N86").Value = -4321
$ws.Range("H89").Value = 71851.82000000001
$ws.Range("I89").Value = 93321.62
$ws.Range("J89").Value = 2075
$ws.Range("K89").Value = 466608.1
$ws.Range("L89").Value = 10375
$ws.Range("M89").Value = -460992.1
$ws.Range("N89").Value = -21607

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9636.059999999999
$ws.Range("I32").Value = 7534.1123
$ws.Range("J32").Value = 26642.727
$ws.Range("K32").Value = 7534.1123
$ws.Range("L32").Value = 26642.727
$ws.Range("M32").Value = -7247.1123
$ws.Range("N32").Value = -27216.727
$ws.Range("H63").Value = 3422.439
$ws.Range("J63").Value = 7997
$ws.Range("L63").Value = 7997
$ws.Range("N63").Value = -9369
$ws.Range("H66").Value = 3422.439
$ws.Range("J66").Value = 7997
$ws.Range("L66").Value = 39985
$ws.Range("N66").Value = -46849
$ws.Range("H97").Value = 968.9167
$ws.Range("I97").Value = 714
$ws.Range("J97").Value = 2753.3333
$ws.Range("K97").Value = 714
$ws.Range("L97").Value = 2753.3333
$ws.Range("M97").Value = -218
$ws.Range("N97").Value = -3745.3333
$ws.Range("H98").Value = 38000
$ws.Range("J98").Value = 38000
$ws.Range("L98").Value = 38000
$ws.Range("N98").Value = -43990

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1650.2424
$ws.Range("I20").Value = 960.0952
$ws.Range("J20").Value = 2858
$ws.Range("K20").Value = 960.0952
$ws.Range("L20").Value = 2858
$ws.Range("M20").Value = -713.0952
$ws.Range("N20").Value = -3352
$ws.Range("H86").Value = 1360.2667
$ws.Range("I86").Value = 1200.3636
$ws.Range("J86").Value = 1800
$ws.Range("K86").Value = 1200.3636
$ws.Range("L86").Value = 1800
$ws.Range("M86").Value = -77.36359999999991
$ws.Range("N86").Value = -4046
$ws.Range("H89").Value = 1360.2667
$ws.Range("I89").Value = 1200.3636
$ws.Range("J89").Value = 1800
$ws.Range("K89").Value = 6001.817999999999
$ws.Range("L89").Value = 9000
$ws.Range("M89").Value = -385.8179999999993
$ws.Range("N89").Value = -20232
$ws.Range("H94").Value = 1737.1818
$ws.Range("I94").Value = 1476.125
$ws.Range("J94").Value = 2433.3333
$ws.Range("K94").Value = 1476.125
$ws.Range("L94").Value = 2433.3333
$ws.Range("M94").Value = -1025.125
$ws.Range("N94").Value = -3335.3333
$ws.Range("H99").Value = 2346.0715
$ws.Range("I99").Value = 648.3333
$ws.Range("J99").Value = 2809.0908
$ws.Range("K99").Value = 648.3333
$ws.Range("L99").Value = 2809.0908
$ws.Range("M99").Value = 849.6667
$ws.Range("N99").Value = -5805.0908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 40757
$ws.Range("J64").Value = 40757
$ws.Range("L64").Value = 40757
$ws.Range("N64").Value = -41253
$ws.Range("H67").Value = 40757
$ws.Range("J67").Value = 40757
$ws.Range("L67").Value = 40757
$ws.Range("N67").Value = -42473
$ws.Range("H68").Value = 51573.75
$ws.Range("J68").Value = 51573.75
$ws.Range("L68").Value = 51573.75
$ws.Range("N68").Value = -53071.75
$ws.Range("H71").Value = 51573.75
$ws.Range("J71").Value = 51573.75
$ws.Range("L71").Value = 154721.25
$ws.Range("N71").Value = -162209.25
$ws.Range("H81").Value = 41666.668
$ws.Range("J81").Value = 41666.668
$ws.Range("L81").Value = 41666.668
$ws.Range("N81").Value = -43662.668
$ws.Range("H84").Value = 41666.668
$ws.Range("J84").Value = 41666.668
$ws.Range("L84").Value = 125000.004
$ws.Range("N84").Value = -134984.004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1383.3334
$ws.Range("J34").Value = 1580
$ws.Range("L34").Value = 4740
$ws.Range("N34").Value = -4908
$ws.Range("H39").Value = 908.53845
$ws.Range("I39").Value = 890.6
$ws.Range("J39").Value = 919.75
$ws.Range("K39").Value = 2671.8
$ws.Range("L39").Value = 2759.25
$ws.Range("M39").Value = -2377.8
$ws.Range("N39").Value = -3347.25
$ws.Range("H44").Value = 266
$ws.Range("I44").Value = 237.5
$ws.Range("J44").Value = 380
$ws.Range("K44").Value = 712.5
$ws.Range("L44").Value = 1140
$ws.Range("M44").Value = -314.5
$ws.Range("N44").Value = -1936
$ws.Range("H55").Value = 896
$ws.Range("I55").Value = 660
$ws.Range("J55").Value = 1250
$ws.Range("K55").Value = 1980
$ws.Range("L55").Value = 3750
$ws.Range("M55").Value = -1803
$ws.Range("N55").Value = -4104
$ws.Range("H131").Value = 1056.4728
$ws.Range("I131").Value = 560
$ws.Range("J131").Value = 1140.9788
$ws.Range("K131").Value = 1680
$ws.Range("L131").Value = 3422.936400000001
$ws.Range("M131").Value = 3360
$ws.Range("N131").Value = -13502.9364

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2470.95
$ws.Range("I80").Value = 2248.4614
$ws.Range("J80").Value = 2884.1428
$ws.Range("K80").Value = 2248.4614
$ws.Range("L80").Value = 2884.1428
$ws.Range("M80").Value = -1250.4614
$ws.Range("N80").Value = -4880.1428
$ws.Range("H83").Value = 2470.95
$ws.Range("I83").Value = 2248.4614
$ws.Range("J83").Value = 2884.1428
$ws.Range("K83").Value = 11242.307
$ws.Range("L83").Value = 14420.714
$ws.Range("M83").Value = -6250.307000000001
$ws.Range("N83").Value = -24404.714
$ws.Range("H123").Value = 15260
$ws.Range("J123").Value = 15260
$ws.Range("L123").Value = 15260
$ws.Range("N123").Value = -20160

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1250731.2
$ws.Range("I46").Value = 475
$ws.Range("J46").Value = 2500987.5
$ws.Range("K46").Value = 475
$ws.Range("L46").Value = 2500987.5
$ws.Range("M46").Value = -287
$ws.Range("N46").Value = -2501363.5
$ws.Range("H68").Value = 2634.25
$ws.Range("I68").Value = 1792.3077
$ws.Range("J68").Value = 3629.2727
$ws.Range("K68").Value = 1792.3077
$ws.Range("L68").Value = 3629.2727
$ws.Range("M68").Value = -1043.3077
$ws.Range("N68").Value = -5127.2727
$ws.Range("H71").Value = 2634.25
$ws.Range("I71").Value = 1792.3077
$ws.Range("J71").Value = 3629.2727
$ws.Range("K71").Value = 8961.538500000001
$ws.Range("L71").Value = 18146.3635
$ws.Range("M71").Value = -5217.538500000001
$ws.Range("N71").Value = -25634.3635

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4363.364
$ws.Range("I62").Value = 3322
$ws.Range("J62").Value = 5231.1665
$ws.Range("K62").Value = 3322
$ws.Range("L62").Value = 5231.1665
$ws.Range("M62").Value = -2698
$ws.Range("N62").Value = -6479.1665
$ws.Range("H65").Value = 4363.364
$ws.Range("I65").Value = 3322
$ws.Range("J65").Value = 5231.1665
$ws.Range("K65").Value = 16610
$ws.Range("L65").Value = 26155.8325
$ws.Range("M65").Value = -13490
$ws.Range("N65").Value = -32395.8325
$ws.Range("H81").Value = 2700.182
$ws.Range("I81").Value = 2522.4443
$ws.Range("J81").Value = 3500
$ws.Range("K81").Value = 5044.8886
$ws.Range("L81").Value = 7000
$ws.Range("M81").Value = -3983.8886
$ws.Range("N81").Value = -9122
$ws.Range("H84").Value = 2700.182
$ws.Range("I84").Value = 2522.4443
$ws.Range("J84").Value = 3500
$ws.Range("K84").Value = 25224.443
$ws.Range("L84").Value = 35000
$ws.Range("M84").Value = -19920.443
$ws.Range("N84").Value = -45608
$ws.Range("H119").Value = 45566.332
$ws.Range("J119").Value = 45566.332
$ws.Range("L119").Value = 45566.332
$ws.Range("N119").Value = -55242.332
